$d = $word.ActiveDocument

# 1) Top-of-document "Date:" field -- scope the find to just that paragraph so
#    the later "Report rendered ..." occurrence of the same date text is untouched.
#    Locate the paragraph dynamically (rather than assuming a fixed index) by
#    scanning for the one that starts with "Date:".
$dateParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Date:")) {
        $dateParaIndex = $i
        break
    }
}
$dateRange = $d.Paragraphs.Item($dateParaIndex).Range
$dateRange.Find.Execute("2017-03-01", $true, $false, $false, $false, $false, $true, 1, $false, "2017-04-24", 2)

# 2) "Report rendered by ..." line (session info)
$d.Content.Find.Execute(
    "Report rendered by koval_000 at 2017-03-01, 09:04 -0500",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Report rendered by koval_000 at 2017-04-24, 11:36 -0400", 2)

# 3) "other attached packages" line [1]
$d.Content.Find.Execute(
    "[1] knitr_1.15.1    forestplot_1.7  checkmate_1.8.2 ggplot2_2.2.1   magrittr_1.5   ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[1] dplyr_0.5.0     forestplot_1.7  checkmate_1.8.2 ggplot2_2.2.1   magrittr_1.5    knitr_1.15.1   ", 2)

# 4) "loaded via a namespace" line [7]
$d.Content.Find.Execute(
    " [7] stringr_1.1.0    plyr_1.8.4       dplyr_0.5.0      tools_3.3.2      DT_0.2           gtable_0.2.0    ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " [7] stringr_1.1.0    plyr_1.8.4       tools_3.3.2      DT_0.2           gtable_0.2.0     DBI_0.5-1       ", 2)

# 5) "loaded via a namespace" line [13]
$d.Content.Find.Execute(
    "[13] plotrix_3.6-4    DBI_0.5-1        htmltools_0.3.5  yaml_2.1.14      lazyeval_0.2.0   assertthat_0.1  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[13] htmltools_0.3.5  yaml_2.1.14      lazyeval_0.2.0   assertthat_0.1   rprojroot_1.2    digest_0.6.12   ", 2)

# 6) "loaded via a namespace" line [19]
$d.Content.Find.Execute(
    "[19] digest_0.6.12    rprojroot_1.2    tibble_1.2       readr_1.0.0      tidyr_0.6.1      htmlwidgets_0.8 ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[19] tibble_1.2       readr_1.0.0      tidyr_0.6.1      htmlwidgets_0.8  evaluate_0.10    haven_1.0.0     ", 2)

# 7) "loaded via a namespace" line [25]
$d.Content.Find.Execute(
    "[25] evaluate_0.10    rmarkdown_1.3    stringi_1.1.2    scales_0.4.1     backports_1.0.5  jsonlite_1.2    ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[25] rmarkdown_1.3    stringi_1.1.2    scales_0.4.1     backports_1.0.5 ", 2)
